$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("K2").Value = 8

# Row 3 changes
$ws.Range("H3").Value = 2.7
$ws.Range("I3").Value = 2.65
$ws.Range("R3").Value = 1.8
$ws.Range("S3").Value = 1.8
